# Add a new client/sale row (row 21 = client #20 "BEBETO SANTOS") to the
# "Cliente" sheet, mirroring the data already present in row 20 (client #19).
#
# Columns: A=# (unchanged, already 20), B=NOME, C=COD - ATIVACAO - MAC,
#          D=DATA INICIO, E=PRAZO ACESSO - DIAS, F=EMAIL, G=CONFIRMADO

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "BEBETO SANTOS"
$ws.Range("C21").Value = "dad3165d0cd51465c2f4943c1053ea42"

$ws.Range("D21").Value = 44852
$ws.Range("D21").NumberFormat = "yyyy-mm-dd"

$ws.Range("E21").Value = 365
$ws.Range("E21").NumberFormat = "#,##0"

$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = "VENDA 17 (18/10)"
